$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 111222590
$ws.Range("I18").Value = 111222590
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 111222590
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -111222306

$ws.Range("H64").Value = 3529.6
$ws.Range("I64").Value = 3249.3333
$ws.Range("J64").Value = 3950
$ws.Range("K64").Value = 3249.3333
$ws.Range("L64").Value = 3950
$ws.Range("M64").Value = -3001.3333
$ws.Range("N64").Value = -4446

$ws.Range("H67").Value = 3529.6
$ws.Range("I67").Value = 3249.3333
$ws.Range("J67").Value = 3950
$ws.Range("K67").Value = 3249.3333
$ws.Range("L67").Value = 3950
$ws.Range("M67").Value = -2391.3333
$ws.Range("N67").Value = -5666

$ws.Range("H132").Value = 2698.6785
$ws.Range("I132").Value = 2945.9565
$ws.Range("J132").Value = 1561.2
$ws.Range("K132").Value = 8837.869499999999
$ws.Range("L132").Value = 4683.6
$ws.Range("M132").Value = -6307.869499999999
$ws.Range("N132").Value = -9743.6

$ws.Range("H137").Value = 42363.48
$ws.Range("I137").Value = 1555.0625
$ws.Range("J137").Value = 114911.78
$ws.Range("K137").Value = 4665.1875
$ws.Range("L137").Value = 344735.34
$ws.Range("M137").Value = -2115.1875
$ws.Range("N137").Value = -349835.34

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12994096
$ws.Range("I32").Value = 17243266
$ws.Range("J32").Value = 22947.947
$ws.Range("K32").Value = 17243266
$ws.Range("L32").Value = 22947.947
$ws.Range("M32").Value = -17242979
$ws.Range("N32").Value = -23521.947

$ws.Range("H63").Value = 3333.3333
$ws.Range("I63").Value = 4000
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 4000
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -3314
$ws.Range("N63").Value = -3372

$ws.Range("H66").Value = 3333.3333
$ws.Range("I66").Value = 4000
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 20000
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -16568
$ws.Range("N66").Value = -16864

$ws.Range("H117").Value = 32979.4
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 32979.4
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 32979.4
$ws.Range("N117").Value = -42157.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2084.4666
$ws.Range("I86").Value = 1855.5834
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1855.5834
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -732.5834
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 2084.4666
$ws.Range("I89").Value = 1855.5834
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 9277.916999999999
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -3661.916999999999
$ws.Range("N89").Value = -26232

$ws.Range("H105").Value = 2884.77
$ws.Range("I105").Value = 2299.7856
$ws.Range("J105").Value = 2980
$ws.Range("K105").Value = 2299.7856
$ws.Range("L105").Value = 2980
$ws.Range("M105").Value = -552.7856000000002
$ws.Range("N105").Value = -6474

$ws.Range("H125").Value = 44350
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 44350
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 44350
$ws.Range("N125").Value = -54190

$ws.Range("H134").Value = 1867.6111
$ws.Range("I134").Value = 1862.1875
$ws.Range("J134").Value = 1911
$ws.Range("K134").Value = 5586.5625
$ws.Range("L134").Value = 5733
$ws.Range("M134").Value = -3051.5625
$ws.Range("N134").Value = -10803

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1709.4572
$ws.Range("I31").Value = 1637
$ws.Range("J31").Value = 1770.4736
$ws.Range("K31").Value = 1637
$ws.Range("L31").Value = 1770.4736
$ws.Range("M31").Value = -1342
$ws.Range("N31").Value = -2360.4736

$ws.Range("H34").Value = 1709.4572
$ws.Range("I34").Value = 1637
$ws.Range("J34").Value = 1770.4736
$ws.Range("K34").Value = 1637
$ws.Range("L34").Value = 1770.4736
$ws.Range("M34").Value = -1435
$ws.Range("N34").Value = -2174.4736

$ws.Range("H86").Value = 3625.875
$ws.Range("I86").Value = 8026.75
$ws.Range("J86").Value = 2158.9167
$ws.Range("K86").Value = 8026.75
$ws.Range("L86").Value = 2158.9167
$ws.Range("M86").Value = -6903.75
$ws.Range("N86").Value = -4404.9167

$ws.Range("H89").Value = 3625.875
$ws.Range("I89").Value = 8026.75
$ws.Range("J89").Value = 2158.9167
$ws.Range("K89").Value = 40133.75
$ws.Range("L89").Value = 10794.5835
$ws.Range("M89").Value = -34517.75
$ws.Range("N89").Value = -22026.5835

$ws.Range("H99").Value = 2205.9355
$ws.Range("I99").Value = 1553.2632
$ws.Range("J99").Value = 3239.3333
$ws.Range("K99").Value = 1553.2632
$ws.Range("L99").Value = 3239.3333
$ws.Range("M99").Value = -55.2632000000001
$ws.Range("N99").Value = -6235.3333

$ws.Range("H126").Value = 2205.9355
$ws.Range("I126").Value = 1553.2632
$ws.Range("J126").Value = 3239.3333
$ws.Range("K126").Value = 4659.7896
$ws.Range("L126").Value = 9717.999899999999
$ws.Range("M126").Value = -2189.7896
$ws.Range("N126").Value = -14657.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 691.5714
$ws.Range("I5").Value = 685.375
$ws.Range("J5").Value = 695.38464
$ws.Range("K5").Value = 2056.125
$ws.Range("L5").Value = 2086.15392
$ws.Range("M5").Value = -1944.125
$ws.Range("N5").Value = -2310.15392

$ws.Range("H68").Value = 1013.0714
$ws.Range("I68").Value = 939.0769
$ws.Range("J68").Value = 1077.2
$ws.Range("K68").Value = 2817.2307
$ws.Range("L68").Value = 3231.6
$ws.Range("M68").Value = -2006.2307
$ws.Range("N68").Value = -4853.6

$ws.Range("H71").Value = 1013.0714
$ws.Range("I71").Value = 939.0769
$ws.Range("J71").Value = 1077.2
$ws.Range("K71").Value = 8451.6921
$ws.Range("L71").Value = 9694.800000000001
$ws.Range("M71").Value = -4395.6921
$ws.Range("N71").Value = -17806.8

$ws.Range("H107").Value = 25641832
$ws.Range("I107").Value = 203.38461
$ws.Range("J107").Value = 38462644
$ws.Range("K107").Value = 610.15383
$ws.Range("L107").Value = 115387932
$ws.Range("M107").Value = 1309.84617
$ws.Range("N107").Value = -115391772

$ws.Range("H135").Value = 691.5714
$ws.Range("I135").Value = 685.375
$ws.Range("J135").Value = 695.38464
$ws.Range("K135").Value = 6168.375
$ws.Range("L135").Value = 6258.46176
$ws.Range("M135").Value = -3633.375
$ws.Range("N135").Value = -11328.46176

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 19900
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 19900
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 19900
$ws.Range("N100").Value = -22064

$ws.Range("H102").Value = 1720.625
$ws.Range("I102").Value = 1789.2632
$ws.Range("J102").Value = 1459.8
$ws.Range("K102").Value = 1789.2632
$ws.Range("L102").Value = 1459.8
$ws.Range("M102").Value = -167.2632000000001
$ws.Range("N102").Value = -4703.8

$ws.Range("H134").Value = 14594
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 14594
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 43782
$ws.Range("N134").Value = -48852

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3742.7188
$ws.Range("I40").Value = 2663.5264
$ws.Range("J40").Value = 5320
$ws.Range("K40").Value = 2663.5264
$ws.Range("L40").Value = 5320
$ws.Range("M40").Value = -2527.5264
$ws.Range("N40").Value = -5592

$ws.Range("H68").Value = 2002192
$ws.Range("I68").Value = 3335000
$ws.Range("J68").Value = 2980
$ws.Range("K68").Value = 3335000
$ws.Range("L68").Value = 2980
$ws.Range("M68").Value = -3334251
$ws.Range("N68").Value = -4478

$ws.Range("H71").Value = 2002192
$ws.Range("I71").Value = 3335000
$ws.Range("J71").Value = 2980
$ws.Range("K71").Value = 16675000
$ws.Range("L71").Value = 14900
$ws.Range("M71").Value = -16671256
$ws.Range("N71").Value = -22388

$ws.Range("H104").Value = 26748
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 26748
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 26748
$ws.Range("N104").Value = -33736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 676.46155
$ws.Range("I107").Value = 674.5
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 2023.5
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = -103.5
$ws.Range("N107").Value = -5940

$ws.Range("H123").Value = 500015000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 500015000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 500015000
$ws.Range("N123").Value = -500024800
